$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the start/end timestamp values (B1/B2); downstream formulas in
# B3 (=B2-B1) and B4 (=B3/POWER(10,9)) recalc automatically.
$ws.Range("B1").Value = 1510872012008000000
$ws.Range("B2").Value = 1510872221873999872

$excel.Calculate()

# Move the active selection to B2
$ws.Range("B2").Select()
